$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.040936523213929
$ws.Range("C2").Value = 0.1440678763283501
$ws.Range("D2").Value = 0.03147556612527325
$ws.Range("E2").Value = 0.0941923848748516
$ws.Range("F2").Value = 3.472502459719266
$ws.Range("J2").Value = 0.1843251799547332
$ws.Range("K2").Value = 1.106647970836605
$ws.Range("M2").Value = 0.3884562439538257
$ws.Range("N2").Value = 3.608084878746737

$ws.Range("B3").Value = 1.002616785891178
$ws.Range("C3").Value = 0.1380182215027759
$ws.Range("D3").Value = 0.03117916304615065
$ws.Range("E3").Value = 0.09372680780220932
$ws.Range("F3").Value = 3.453374209637289
$ws.Range("J3").Value = 0.1836150428679844
$ws.Range("K3").Value = 1.064745862659692
$ws.Range("M3").Value = 0.3789740781784872
$ws.Range("N3").Value = 3.61270804478049

$ws.Range("B4").Value = 0.9797222471354701
$ws.Range("C4").Value = 0.1343971755217552
$ws.Range("D4").Value = 0.03101622485688438
$ws.Range("E4").Value = 0.09348286438027564
$ws.Range("F4").Value = 3.442993684828451
$ws.Range("J4").Value = 0.1832525573495829
$ws.Range("K4").Value = 1.03970120394527
$ws.Range("M4").Value = 0.3733789150369304
$ws.Range("N4").Value = 3.616319524722684

$ws.Range("B5").Value = 0.9705518261576742
$ws.Range("C5").Value = 0.1329450154798906
$ws.Range("D5").Value = 0.03095464164490735
$ws.Range("E5").Value = 0.09339401245362389
$ws.Range("F5").Value = 3.439106426931914
$ws.Range("J5").Value = 0.183123337891999
$ws.Range("K5").Value = 1.029666963550881
$ws.Range("M5").Value = 0.3711558974214597
$ws.Range("N5").Value = 3.617985388974787

$ws.Range("B6").Value = 0.9690387021583149
$ws.Range("C6").Value = 0.1327052998244795
$ws.Range("D6").Value = 0.03094470734889398
$ws.Range("E6").Value = 0.09337989672773617
$ws.Range("F6").Value = 3.438481655583601
$ws.Range("J6").Value = 0.1831029983537178
$ws.Range("K6").Value = 1.02801114826778
$ws.Range("M6").Value = 0.3707902130347307
$ws.Range("N6").Value = 3.618273725887036

$ws.Range("B7").Value = 0.9795979268280064
$ws.Range("C7").Value = 0.1343774963276871
$ws.Range("D7").Value = 0.03101537479412286
$ws.Range("E7").Value = 0.09348162332676857
$ws.Range("F7").Value = 3.44293987180221
$ws.Range("J7").Value = 0.1832507397546195
$ws.Range("K7").Value = 1.03956518386255
$ws.Range("M7").Value = 0.3733487036287642
$ws.Range("N7").Value = 3.616341205275674

$ws.Range("B8").Value = 1.027592216632002
$ws.Range("C8").Value = 0.1419625192838225
$ws.Range("D8").Value = 0.03136942572517398
$ws.Range("E8").Value = 0.0940231614929985
$ws.Range("F8").Value = 3.465623759463455
$ws.Range("J8").Value = 0.1840650626496227
$ws.Range("K8").Value = 1.092058152974886
$ws.Range("M8").Value = 0.3851396870669319
$ws.Range("N8").Value = 3.609518432469969

$ws.Range("B9").Value = 1.126750821029333
$ws.Range("C9").Value = 0.1575819356266095
$ws.Range("D9").Value = 0.03221394087337615
$ws.Range("E9").Value = 0.09541725817980051
$ws.Range("F9").Value = 3.520945794967929
$ws.Range("J9").Value = 0.1862455951661772
$ws.Range("K9").Value = 1.200434798318298
$ws.Range("M9").Value = 0.410064408209486
$ws.Range("N9").Value = 3.602280622564493

$ws.Range("B10").Value = 1.202702450768072
$ws.Range("C10").Value = 0.169518252604405
$ws.Range("D10").Value = 0.03292481664953328
$ws.Range("E10").Value = 0.09664356862871415
$ws.Range("F10").Value = 3.568226683539621
$ws.Range("J10").Value = 0.1882040475264972
$ws.Range("K10").Value = 1.283406248854845
$ws.Range("M10").Value = 0.429481104229275
$ws.Range("N10").Value = 3.600722327107007

$ws.Range("B11").Value = 1.237934210162962
$ws.Range("C11").Value = 0.1750499290519372
$ws.Range("D11").Value = 0.0332676081329808
$ws.Range("E11").Value = 0.09724525919538962
$ws.Range("F11").Value = 3.591183626566959
$ws.Range("J11").Value = 0.1891725567466622
$ws.Range("K11").Value = 1.321886528387779
$ws.Range("M11").Value = 0.4385554914805141
$ws.Range("N11").Value = 3.60083321758195

$ws.Range("B12").Value = 1.25137380740324
$ws.Range("C12").Value = 0.1771593605082558
$ws.Range("D12").Value = 0.03340018388582422
$ws.Range("E12").Value = 0.09747939734690547
$ws.Range("F12").Value = 3.600085530014979
$ws.Range("J12").Value = 0.1895504710762381
$ws.Range("K12").Value = 1.336564280185712
$ws.Range("M12").Value = 0.4420265348665495
$ws.Range("N12").Value = 3.600993353464816

$ws.Range("B13").Value = 1.248474984104007
$ws.Range("C13").Value = 0.1767044015520582
$ws.Range("D13").Value = 0.03337150858120452
$ws.Range("E13").Value = 0.09742869191006065
$ws.Range("F13").Value = 3.598159063141878
$ws.Range("J13").Value = 0.1894685840391688
$ws.Range("K13").Value = 1.333398440637581
$ws.Range("M13").Value = 0.4412774367781864
$ws.Range("N13").Value = 3.600953606608172

$ws.Range("B14").Value = 1.23903792695512
$ws.Range("C14").Value = 0.1752231781253499
$ws.Range("D14").Value = 0.03327845987533351
$ws.Range("E14").Value = 0.09726439587199209
$ws.Range("F14").Value = 3.591911808989948
$ws.Range("J14").Value = 0.1892034242906604
$ws.Range("K14").Value = 1.32309194791182
$ws.Range("M14").Value = 0.4388403593288075
$ws.Range("N14").Value = 3.600844022754032

$ws.Range("B15").Value = 1.233270236943156
$ws.Range("C15").Value = 0.1743178031856871
$ws.Range("D15").Value = 0.03322182465642953
$ws.Range("E15").Value = 0.09716457871024531
$ws.Range("F15").Value = 3.588112361220396
$ws.Range("J15").Value = 0.1890424598412537
$ws.Range("K15").Value = 1.316792749018617
$ws.Range("M15").Value = 0.4373521078081879
$ws.Range("N15").Value = 3.600792293197784

$ws.Range("B16").Value = 1.200413689924744
$ws.Range("C16").Value = 0.1691587985933722
$ws.Range("D16").Value = 0.0329028030249674
$ws.Range("E16").Value = 0.09660512774502195
$ws.Range("F16").Value = 3.566755568557255
$ws.Range("J16").Value = 0.1881423148480081
$ws.Range("K16").Value = 1.28090630629697
$ws.Range("M16").Value = 0.4288929359153855
$ws.Range("N16").Value = 3.600731598292739

$ws.Range("B17").Value = 1.180431797296251
$ws.Range("C17").Value = 0.1660200348086676
$ws.Range("D17").Value = 0.03271204853791687
$ws.Range("E17").Value = 0.09627314108558949
$ws.Range("F17").Value = 3.554025104543427
$ws.Range("J17").Value = 0.1876099815406249
$ws.Range("K17").Value = 1.259079855439523
$ws.Range("M17").Value = 0.4237654121635259
$ws.Range("N17").Value = 3.600904511160948

$ws.Range("B18").Value = 1.169002827821117
$ws.Range("C18").Value = 0.1642242772688292
$ws.Range("D18").Value = 0.03260415998040145
$ws.Range("E18").Value = 0.09608631790550604
$ws.Range("F18").Value = 3.546839211979986
$ws.Range("J18").Value = 0.1873111013668733
$ws.Range("K18").Value = 1.246595132885773
$ws.Range("M18").Value = 0.4208389384247511
$ws.Range("N18").Value = 3.601081110432915

$ws.Range("B19").Value = 1.165144177003612
$ws.Range("C19").Value = 0.1636179072934567
$ws.Range("D19").Value = 0.03256794550475206
$ws.Range("E19").Value = 0.09602377194917366
$ws.Range("F19").Value = 3.544429598773206
$ws.Range("J19").Value = 0.1872111600612172
$ws.Range("K19").Value = 1.242379910925905
$ws.Range("M19").Value = 0.4198519900443287
$ws.Range("N19").Value = 3.601154145594791

$ws.Range("B20").Value = 1.182552269947962
$ws.Range("C20").Value = 0.1663531699589385
$ws.Range("D20").Value = 0.03273216561607484
$ws.Range("E20").Value = 0.09630805461943126
$ws.Range("F20").Value = 3.555366172084817
$ws.Range("J20").Value = 0.1876658933935715
$ws.Range("K20").Value = 1.261396146965382
$ws.Range("M20").Value = 0.4243088919196936
$ws.Range("N20").Value = 3.600878118300614

$ws.Range("B21").Value = 1.241807155747722
$ws.Range("C21").Value = 0.1756578496876102
$ws.Range("D21").Value = 0.03330571559491347
$ws.Range("E21").Value = 0.09731248296394668
$ws.Range("F21").Value = 3.593741115018247
$ws.Range("J21").Value = 0.189281005212564
$ws.Range("K21").Value = 1.326116332984583
$ws.Range("M21").Value = 0.4395552442122437
$ws.Range("N21").Value = 3.600873001674032

$ws.Range("B22").Value = 1.281105503327865
$ws.Range("C22").Value = 0.1818247445086456
$ws.Range("D22").Value = 0.03369668745627052
$ws.Range("E22").Value = 0.09800559804652664
$ws.Range("F22").Value = 3.620037426138254
$ws.Range("J22").Value = 0.1904016298296582
$ws.Range("K22").Value = 1.369033298086919
$ws.Range("M22").Value = 0.4497222945052783
$ws.Range("N22").Value = 3.601558417073136

$ws.Range("B23").Value = 1.260078852868389
$ws.Range("C23").Value = 0.1785254875138378
$ws.Range("D23").Value = 0.03348655067213002
$ws.Range("E23").Value = 0.0976323188778494
$ws.Range("F23").Value = 3.605891219436529
$ws.Range("J23").Value = 0.1897975779281111
$ws.Range("K23").Value = 1.346071021915606
$ws.Range("M23").Value = 0.4442773952003236
$ws.Range("N23").Value = 3.601129488439057

$ws.Range("B24").Value = 1.181593420852437
$ws.Range("C24").Value = 0.1662025322792715
$ws.Range("D24").Value = 0.03272306514259782
$ws.Range("E24").Value = 0.09629225762494897
$ws.Range("F24").Value = 3.554759461120156
$ws.Range("J24").Value = 0.1876405933264635
$ws.Range("K24").Value = 1.260348753498647
$ws.Range("M24").Value = 0.4240631180418006
$ws.Range("N24").Value = 3.600889810081895

$ws.Range("B25").Value = 1.099383343701874
$ws.Range("C25").Value = 0.1532760300545135
$ws.Range("D25").Value = 0.03196951659402458
$ws.Range("E25").Value = 0.09500460391909726
$ws.Range("F25").Value = 3.504816464229137
$ws.Range("J25").Value = 0.1855931587326296
$ws.Range("K25").Value = 1.170530588199767
$ws.Range("M25").Value = 0.4031280583612897
$ws.Range("N25").Value = 3.603579582798574

